# Update latest output (run 61)
$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates ---
$wsSchedule.Range("E4").Value = 194.72505975
$wsSchedule.Range("F4").Value = 7.359223724489796

$wsSchedule.Range("A5").Value = 46040.29166666666
$wsSchedule.Range("B5").Value = 46040.79166666666
$wsSchedule.Range("E5").Value = 95.10621899999998
$wsSchedule.Range("F5").Value = 2.096697949735449

# --- Detailed sheet updates ---
$wsDetailed.Range("B31").Value = 36.06031
$wsDetailed.Range("B32").Value = 36.06033

$wsDetailed.Range("B33").Value = -9.99
$wsDetailed.Range("C33").Value = "historical"

$wsDetailed.Range("B34").Value = -6
$wsDetailed.Range("C34").Value = "historical"

$wsDetailed.Range("B35").Value = 2.55649
$wsDetailed.Range("B36").Value = -1.77936
$wsDetailed.Range("B37").Value = -0.289
$wsDetailed.Range("B38").Value = -2.33161
$wsDetailed.Range("B39").Value = -2.08153
$wsDetailed.Range("B40").Value = 13.59536
$wsDetailed.Range("B41").Value = 46.95824
$wsDetailed.Range("B42").Value = 47.19413
$wsDetailed.Range("B44").Value = 47.00942
$wsDetailed.Range("B45").Value = 46.87831
$wsDetailed.Range("B47").Value = 47.87261
$wsDetailed.Range("B49").Value = 56.78
$wsDetailed.Range("B50").Value = 48.09522
$wsDetailed.Range("B52").Value = 56.98
$wsDetailed.Range("B55").Value = 56.97996

$wsDetailed.Range("B64").Value = 24.30295
$wsDetailed.Range("E64").Value = "ON"

$wsDetailed.Range("B65").Value = 0.7
$wsDetailed.Range("B66").Value = 5.738
$wsDetailed.Range("B67").Value = 27.1824
$wsDetailed.Range("B68").Value = 14.8839
$wsDetailed.Range("B69").Value = 0.7
$wsDetailed.Range("B70").Value = 25.82334
$wsDetailed.Range("B71").Value = 22.07
$wsDetailed.Range("B72").Value = 22.07
$wsDetailed.Range("B73").Value = 25.80412
$wsDetailed.Range("B74").Value = 22.07
$wsDetailed.Range("B75").Value = -0.93264
$wsDetailed.Range("B76").Value = -5.45227
$wsDetailed.Range("B77").Value = -5.5808

$wsDetailed.Range("B80").Value = -12.01
$wsDetailed.Range("B81").Value = -11.01
$wsDetailed.Range("B82").Value = -7.17238
$wsDetailed.Range("B83").Value = -8.324020000000001
$wsDetailed.Range("B84").Value = -6.1466
$wsDetailed.Range("B85").Value = -6.5703
$wsDetailed.Range("B86").Value = -3.10096

$wsDetailed.Range("B88").Value = 36.0601
$wsDetailed.Range("E88").Value = "OFF"

$wsDetailed.Range("B89").Value = 45.86636
$wsDetailed.Range("B90").Value = 54.96565
$wsDetailed.Range("B91").Value = 47.16323
$wsDetailed.Range("B92").Value = 46.5469
$wsDetailed.Range("B94").Value = 55.29915
